$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1905.3334
$ws.Range("I38").Value = 74.833336
$ws.Range("J38").Value = 5566.3335
$ws.Range("K38").Value = 224.500008
$ws.Range("L38").Value = 16699.0005
$ws.Range("M38").Value = 147.499992
$ws.Range("N38").Value = -17443.0005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 274.6
$ws.Range("I39").Value = 249.55556
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 748.66668
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -452.66668
$ws.Range("N39").Value = -2092

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1994
$ws.Range("J43").Value = 1994
$ws.Range("L43").Value = 1994
$ws.Range("N43").Value = -2132

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4867.4546
$ws.Range("I62").Value = 4427.1665
$ws.Range("J62").Value = 5395.8
$ws.Range("K62").Value = 4427.1665
$ws.Range("L62").Value = 5395.8
$ws.Range("M62").Value = -3803.1665
$ws.Range("N62").Value = -6643.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4867.4546
$ws.Range("I65").Value = 4427.1665
$ws.Range("J65").Value = 5395.8
$ws.Range("K65").Value = 22135.8325
$ws.Range("L65").Value = 26979
$ws.Range("M65").Value = -19015.8325
$ws.Range("N65").Value = -33219

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 672
$ws.Range("I111").Value = 599.2857
$ws.Range("K111").Value = 1797.8571
$ws.Range("M111").Value = 1269.1429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1957.6
$ws.Range("I132").Value = 1980.8334
$ws.Range("K132").Value = 5942.5002
$ws.Range("M132").Value = -3412.5002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3584.4138
$ws.Range("I138").Value = 2196
$ws.Range("J138").Value = 3806.56
$ws.Range("K138").Value = 6588
$ws.Range("L138").Value = 11419.68
$ws.Range("M138").Value = -1448
$ws.Range("N138").Value = -21699.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3255.182
$ws.Range("I74").Value = 3333.7
$ws.Range("K74").Value = 3333.7
$ws.Range("M74").Value = -2459.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3255.182
$ws.Range("I77").Value = 3333.7
$ws.Range("K77").Value = 16668.5
$ws.Range("M77").Value = -12300.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 65429.25
$ws.Range("J130").Value = 65429.25
$ws.Range("L130").Value = 65429.25
$ws.Range("N130").Value = -75469.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1407.3846
$ws.Range("I132").Value = 1290.1428
$ws.Range("K132").Value = 3870.4284
$ws.Range("M132").Value = -1340.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4692.4443
$ws.Range("I20").Value = 3748.8
$ws.Range("K20").Value = 3748.8
$ws.Range("M20").Value = -3501.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 1100
$ws.Range("J30").Value = 1100
$ws.Range("L30").Value = 1100
$ws.Range("N30").Value = -1350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2999.5557
$ws.Range("I94").Value = 2999.5557
$ws.Range("K94").Value = 2999.5557
$ws.Range("M94").Value = -2548.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3880
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 3600
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 3600
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -4848

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3880
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 3600
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 18000
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -24240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H130").Value = 55000
$ws.Range("J130").Value = 55000
$ws.Range("L130").Value = 55000
$ws.Range("N130").Value = -65040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2544.3333
$ws.Range("I132").Value = 2424.875
$ws.Range("K132").Value = 7274.625
$ws.Range("M132").Value = -4744.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 316
$ws.Range("I36").Value = 324
$ws.Range("J36").Value = 300
$ws.Range("K36").Value = 972
$ws.Range("L36").Value = 900
$ws.Range("M36").Value = -803
$ws.Range("N36").Value = -1238

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1030
$ws.Range("I47").Value = 545
$ws.Range("K47").Value = 1635
$ws.Range("M47").Value = -1204

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1030
$ws.Range("I131").Value = 1030
$ws.Range("K131").Value = 3090
$ws.Range("M131").Value = 1950

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("N26").Value = -10560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 10000
$ws.Range("J50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("N50").Value = -10996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1168.75
$ws.Range("I97").Value = 1025
$ws.Range("J97").Value = 1600
$ws.Range("K97").Value = 1025
$ws.Range("L97").Value = 1600
$ws.Range("M97").Value = -529
$ws.Range("N97").Value = -2592

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4032.75
$ws.Range("J122").Value = 4881
$ws.Range("L122").Value = 14643
$ws.Range("N122").Value = -19543

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 45430
$ws.Range("J128").Value = 45430
$ws.Range("L128").Value = 45430
$ws.Range("N128").Value = -55390

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5594.9546
$ws.Range("I132").Value = 4354.8184
$ws.Range("J132").Value = 6835.091
$ws.Range("K132").Value = 13064.4552
$ws.Range("L132").Value = 20505.273
$ws.Range("M132").Value = -10534.4552
$ws.Range("N132").Value = -25565.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3019
$ws.Range("I136").Value = 2740.2856
$ws.Range("J136").Value = 3994.5
$ws.Range("K136").Value = 8220.856800000001
$ws.Range("L136").Value = 11983.5
$ws.Range("M136").Value = -5670.856800000001
$ws.Range("N136").Value = -17083.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6470.8125
$ws.Range("I136").Value = 6441
$ws.Range("J136").Value = 6600
$ws.Range("K136").Value = 19323
$ws.Range("L136").Value = 19800
$ws.Range("M136").Value = -16773
$ws.Range("N136").Value = -24900
